{"js": "// The document repeats a \"cost profile\" chart page 4 times (Picture 1..4);\n// every copy after the first is preceded by its own section-break-only\n// paragraph so that picture could carry its own page size/orientation. This\n// edit keeps only the first chart page and removes the other three picture\n// paragraphs together with the section-break paragraphs that separated them,\n// then fixes the remaining (now last) section's page size so it reads\n// 15840 x 12240 landscape instead of 12240 x 15840 landscape.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph holding the first inline picture \u2014 everything after\n// it is a repeat of the same chart page and gets dropped.\nlet firstPictureIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const pics = paragraphs.items[i].inlinePictures;\n  pics.load(\"items\");\n  await context.sync();\n  if (pics.items.length > 0) {\n    firstPictureIndex = i;\n    break;\n  }\n}\n\n// Paragraph layout (0-based) before the edit confirms this lands on index 2:\n//   0: empty paragraph\n//   1: paragraph holding the first section break (portrait)\n//   2: paragraph with Picture 1   <- keep everything up to here\n//   3: paragraph holding a section break (landscape)   <- delete from here\n//   4: paragraph with Picture 2                        <- delete\n//   5: paragraph holding a section break (landscape)    <- delete\n//   6: paragraph with Picture 3                         <- delete\n//   7: paragraph holding a section break (landscape)     <- delete\n//   8: paragraph with Picture 4                          <- delete\n// Delete from the end backward so earlier indices stay valid.\nif (firstPictureIndex >= 0) {\n  for (let i = paragraphs.items.length - 1; i > firstPictureIndex; i--) {\n    paragraphs.items[i].delete();\n  }\n  await context.sync();\n}\n\n// The final body section now inherits the page size that used to belong to\n// the removed tail; correct it to 15840 x 12240 (landscape) \u2014 i.e. swap the\n// existing 12240 x 15840 width/height (values are in points: 1 twip = 1/20\n// pt, so 15840 twips = 792 pt and 12240 twips = 612 pt).\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst lastSection = sections.items[sections.items.length - 1];\nconst pageSetup = lastSection.pageSetup;\npageSetup.pageWidth = 792;\npageSetup.pageHeight = 612;\nawait context.sync();\n", "ps1": "# The document repeats a \"cost profile\" chart page 4 times (Picture 1..4);\n# every copy after the first is preceded by its own section-break-only\n# paragraph so that picture could carry its own page size/orientation. This\n# edit keeps only the first chart page and removes the other three picture\n# paragraphs together with the section-break paragraphs that separated them,\n# then fixes the remaining (now last) section's page size so it reads\n# 15840 x 12240 landscape instead of 12240 x 15840 landscape.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph holding the first inline shape (picture) -- everything\n# after it is a repeat of the same chart page and gets dropped.\n$firstPictureIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  if ($d.Paragraphs.Item($i).Range.InlineShapes.Count -gt 0) {\n    $firstPictureIndex = $i\n    break\n  }\n}\n\n# Paragraph layout (1-based) before the edit confirms this lands on index 3:\n#   1: empty paragraph\n#   2: paragraph holding the first section break (portrait)\n#   3: paragraph with Picture 1   <- keep everything up to here\n#   4: paragraph holding a section break (landscape)   <- delete from here\n#   5: paragraph with Picture 2                        <- delete\n#   6: paragraph holding a section break (landscape)    <- delete\n#   7: paragraph with Picture 3                         <- delete\n#   8: paragraph holding a section break (landscape)     <- delete\n#   9: paragraph with Picture 4                          <- delete\n# Delete from the end backward so earlier indices stay valid.\nif ($firstPictureIndex -gt 0) {\n  for ($i = $d.Paragraphs.Count; $i -ge $firstPictureIndex + 1; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n  }\n}\n\n# The final section now inherits the page size that used to belong to the\n# removed tail; correct it to 15840 x 12240 (landscape) -- i.e. swap the\n# existing 12240 x 15840 width/height (values are in points: 1 twip = 1/20\n# pt, so 15840 twips = 792 pt and 12240 twips = 612 pt).\n$lastSection = $d.Sections.Item($d.Sections.Count)\n$pageSetup = $lastSection.PageSetup\n$pageSetup.PageWidth = 792\n$pageSetup.PageHeight = 612\n"}
